$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 424, shifting existing rows 424..514 down to 425..515.
$ws.Rows("424:424").Insert()

# Populate the newly inserted row 424 with the new data record.
$ws.Range("A424").Value = 7
$ws.Range("B424").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C424").Value = "Ñuble"
$ws.Range("D424").Value = 45173
$ws.Range("E424").Value = 16
$ws.Range("F424").Value = 100114013
$ws.Range("G424").Value = "Zanahoria"
$ws.Range("H424").Value = "Sin especificar"
$ws.Range("I424").Value = "Primera"
$ws.Range("J424").Value = 150
$ws.Range("K424").Value = 6000
$ws.Range("L424").Value = 6000
$ws.Range("M424").Value = 6000
$ws.Range("N424").Value = "$/saco 20 kilos"
$ws.Range("O424").Value = "Provincia de Diguillín"
$ws.Range("P424").Value = 300
$ws.Range("Q424").Value = 20
$ws.Range("R424").Value = "Hortaliza"
